# Auto-generated Excel COM-interop script applying the numeric restatements
# described by the commit "chore: update Sheets via scheduled runner".
# Each (sheet, row) block updates the H-N "leve profit" computed columns to the
# freshly recalculated market-board figures; some rows also gain or lose a
# trailing M/N cell because the new NQ/HQ comparison no longer (or now does)
# produce a profit figure for that branch.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1226  # was 1433.8
$ws.Range("I28").Value = 1320.1818  # was 1572
$ws.Range("K28").Value = 1320.1818  # was 1572
$ws.Range("M28").Value = -835.1818000000001  # was -1087

$ws.Range("H55").Value = 298.1111  # was 288.5
$ws.Range("J55").Value = 384.6  # was 387
$ws.Range("L55").Value = 384.6  # was 387
$ws.Range("N55").Value = -812.6  # was -815

$ws.Range("H62").Value = 44630.8  # was 39051.17
$ws.Range("I62").Value = 81092.53999999999  # was 75584.42999999999
$ws.Range("J62").Value = 5130.5835  # was 4953.467
$ws.Range("K62").Value = 81092.53999999999  # was 75584.42999999999
$ws.Range("L62").Value = 5130.5835  # was 4953.467
$ws.Range("M62").Value = -80468.53999999999  # was -74960.42999999999
$ws.Range("N62").Value = -6378.5835  # was -6201.467

$ws.Range("H65").Value = 44630.8  # was 39051.17
$ws.Range("I65").Value = 81092.53999999999  # was 75584.42999999999
$ws.Range("J65").Value = 5130.5835  # was 4953.467
$ws.Range("K65").Value = 405462.7  # was 377922.15
$ws.Range("L65").Value = 25652.9175  # was 24767.335
$ws.Range("M65").Value = -402342.7  # was -374802.15
$ws.Range("N65").Value = -31892.9175  # was -31007.335

$ws.Range("H98").Value = 1875.1333  # was 1875.2
$ws.Range("I98").Value = 1875.1333  # was 1937.7858
$ws.Range("J98").Value = 0  # was 999
$ws.Range("K98").Value = 1875.1333  # was 1937.7858
$ws.Range("L98").Value = 0  # was 999
$ws.Range("M98").Value = -377.1333  # was -439.7858000000001
$ws.Range("N98").ClearContents()  # was -3995

$ws.Range("H122").Value = 1875.1333  # was 1875.2
$ws.Range("I122").Value = 1875.1333  # was 1937.7858
$ws.Range("J122").Value = 0  # was 999
$ws.Range("K122").Value = 5625.3999  # was 5813.357400000001
$ws.Range("L122").Value = 0  # was 2997
$ws.Range("M122").Value = -3175.3999  # was -3363.357400000001
$ws.Range("N122").ClearContents()  # was -7897

$ws.Range("H129").Value = 1061.1052  # was 1061.2632
$ws.Range("I129").Value = 473.17648  # was 473.35294
$ws.Range("K129").Value = 1419.52944  # was 1420.05882
$ws.Range("M129").Value = 3580.47056  # was 3579.94118

$ws.Range("H132").Value = 1789.3334  # was 1808.65
$ws.Range("I132").Value = 1768.6842  # was 1789
$ws.Range("K132").Value = 5306.0526  # was 5367
$ws.Range("M132").Value = -2776.0526  # was -2837

$ws.Range("H137").Value = 3115  # was 3094.4429
$ws.Range("I137").Value = 2350.75  # was 2360.0962
$ws.Range("J137").Value = 5452.706  # was 5215.8887
$ws.Range("K137").Value = 7052.25  # was 7080.2886
$ws.Range("L137").Value = 16358.118  # was 15647.6661
$ws.Range("M137").Value = -4502.25  # was -4530.2886
$ws.Range("N137").Value = -21458.118  # was -20747.6661

$ws.Range("H138").Value = 2239.1035  # was 2251.0327
$ws.Range("J138").Value = 2690.6858  # was 2674.1843
$ws.Range("L138").Value = 8072.057400000001  # was 8022.5529
$ws.Range("N138").Value = -18352.0574  # was -18302.5529

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3690.5833  # was 3821.5454
$ws.Range("I2").Value = 3880.375  # was 4113.2856
$ws.Range("K2").Value = 3880.375  # was 4113.2856
$ws.Range("M2").Value = -3767.375  # was -4000.2856

$ws.Range("H32").Value = 24683.457  # was 26247.4
$ws.Range("I32").Value = 4978.569  # was 5282.4585
$ws.Range("J32").Value = 192175  # was 170007
$ws.Range("K32").Value = 4978.569  # was 5282.4585
$ws.Range("L32").Value = 192175  # was 170007
$ws.Range("M32").Value = -4691.569  # was -4995.4585
$ws.Range("N32").Value = -192749  # was -170581

$ws.Range("H45").Value = 563643.4  # was 507436.1
$ws.Range("I45").Value = 1011956.4  # was 843558.75
$ws.Range("K45").Value = 1011956.4  # was 843558.75
$ws.Range("M45").Value = -1011579.4  # was -843181.75

$ws.Range("H74").Value = 1620.7142  # was 1321.3334
$ws.Range("I74").Value = 1432.9333  # was 1177.8276
$ws.Range("J74").Value = 2090.1667  # was 1915.8572
$ws.Range("K74").Value = 1432.9333  # was 1177.8276
$ws.Range("L74").Value = 2090.1667  # was 1915.8572
$ws.Range("M74").Value = -558.9332999999999  # was -303.8276000000001
$ws.Range("N74").Value = -3838.1667  # was -3663.8572

$ws.Range("H77").Value = 1620.7142  # was 1321.3334
$ws.Range("I77").Value = 1432.9333  # was 1177.8276
$ws.Range("J77").Value = 2090.1667  # was 1915.8572
$ws.Range("K77").Value = 7164.666499999999  # was 5889.138000000001
$ws.Range("L77").Value = 10450.8335  # was 9579.286
$ws.Range("M77").Value = -2796.666499999999  # was -1521.138000000001
$ws.Range("N77").Value = -19186.8335  # was -18315.286

$ws.Range("H102").Value = 3167.0833  # was 3300.4546
$ws.Range("I102").Value = 2798.6  # was 2920.6667
$ws.Range("K102").Value = 2798.6  # was 2920.6667
$ws.Range("M102").Value = -1176.6  # was -1298.6667

$ws.Range("H110").Value = 1724.7  # was 1768.375
$ws.Range("I110").Value = 1664.1428  # was 1709.8
$ws.Range("K110").Value = 1664.1428  # was 1709.8
$ws.Range("M110").Value = 380.8571999999999  # was 335.2

$ws.Range("H116").Value = 3690.5833  # was 3821.5454
$ws.Range("I116").Value = 3880.375  # was 4113.2856
$ws.Range("K116").Value = 3880.375  # was 4113.2856
$ws.Range("M116").Value = -1586.375  # was -1819.2856

$ws.Range("H122").Value = 2791.4375  # was 2803.625
$ws.Range("I122").Value = 2827.5454  # was 2845.2727
$ws.Range("K122").Value = 8482.636200000001  # was 8535.8181
$ws.Range("M122").Value = -6032.636200000001  # was -6085.8181

$ws.Range("H132").Value = 2329.3333  # was 1996.5
$ws.Range("I132").Value = 2496.5  # was 1998
$ws.Range("K132").Value = 7489.5  # was 5994
$ws.Range("M132").Value = -4959.5  # was -3464

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3690.5833  # was 3821.5454
$ws.Range("I3").Value = 3880.375  # was 4113.2856
$ws.Range("K3").Value = 3880.375  # was 4113.2856
$ws.Range("M3").Value = -3766.375  # was -3999.2856

$ws.Range("H42").Value = 123999  # was 0
$ws.Range("J42").Value = 123999  # was 0
$ws.Range("L42").Value = 123999  # was 0
$ws.Range("N42").Value = -124655  # newly populated cell

$ws.Range("H82").Value = 13625.417  # was 15248.929
$ws.Range("J82").Value = 25000  # was 24996.666
$ws.Range("L82").Value = 25000  # was 24996.666
$ws.Range("N82").Value = -25766  # was -25762.666

$ws.Range("H85").Value = 13625.417  # was 15248.929
$ws.Range("J85").Value = 25000  # was 24996.666
$ws.Range("L85").Value = 25000  # was 24996.666
$ws.Range("N85").Value = -27652  # was -27648.666

$ws.Range("H94").Value = 1684.25  # was 1413.1818
$ws.Range("I94").Value = 1656.5  # was 1226.5
$ws.Range("J94").Value = 1712  # was 1911
$ws.Range("K94").Value = 1656.5  # was 1226.5
$ws.Range("L94").Value = 1712  # was 1911
$ws.Range("M94").Value = -1205.5  # was -775.5
$ws.Range("N94").Value = -2614  # was -2813

$ws.Range("H102").Value = 4650  # was 13304
$ws.Range("J102").Value = 0  # was 30612
$ws.Range("L102").Value = 0  # was 30612
$ws.Range("N102").ClearContents()  # was -37102

$ws.Range("H105").Value = 1351.8889  # was 1565
$ws.Range("I105").Value = 1036.7142  # was 1157.3334
$ws.Range("J105").Value = 2455  # was 4011
$ws.Range("K105").Value = 1036.7142  # was 1157.3334
$ws.Range("L105").Value = 2455  # was 4011
$ws.Range("M105").Value = 710.2858000000001  # was 589.6666
$ws.Range("N105").Value = -5949  # was -7505

$ws.Range("H128").Value = 14654.777  # was 12624
$ws.Range("I128").Value = 14654.777  # was 12624
$ws.Range("K128").Value = 43964.331  # was 37872
$ws.Range("M128").Value = -41474.331  # was -35382

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 46614.176  # was 32065.266
$ws.Range("I31").Value = 335351.34  # was 73153.36
$ws.Range("K31").Value = 335351.34  # was 73153.36
$ws.Range("M31").Value = -335056.34  # was -72858.36

$ws.Range("H34").Value = 46614.176  # was 32065.266
$ws.Range("I34").Value = 335351.34  # was 73153.36
$ws.Range("K34").Value = 335351.34  # was 73153.36
$ws.Range("M34").Value = -335149.34  # was -72951.36

$ws.Range("H122").Value = 37627.93  # was 36399.367
$ws.Range("I122").Value = 47636.637  # was 45599
$ws.Range("K122").Value = 142909.911  # was 136797
$ws.Range("M122").Value = -140459.911  # was -134347

$ws.Range("H132").Value = 2451.56  # was 2476.7346
$ws.Range("I132").Value = 2376.2683  # was 2405.225
$ws.Range("K132").Value = 7128.804900000001  # was 7215.674999999999
$ws.Range("M132").Value = -4598.804900000001  # was -4685.674999999999

$ws.Range("H134").Value = 2576.907  # was 2606.6428
$ws.Range("I134").Value = 2488.3057  # was 2521.457
$ws.Range("K134").Value = 7464.9171  # was 7564.370999999999
$ws.Range("M134").Value = -4929.9171  # was -5029.370999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1684.2222  # was 1524.8
$ws.Range("J2").Value = 3778.5  # was 3040.8
$ws.Range("L2").Value = 22671  # was 18244.8
$ws.Range("N2").Value = -22897  # was -18470.8

$ws.Range("H45").Value = 10750  # was 12333.333
$ws.Range("J45").Value = 11041.667  # was 13000
$ws.Range("L45").Value = 33125.001  # was 39000
$ws.Range("N45").Value = -34189.001  # was -40064

$ws.Range("H97").Value = 2475.4  # was 2540.6667
$ws.Range("J97").Value = 3097  # was 3500
$ws.Range("L97").Value = 9291  # was 10500
$ws.Range("N97").Value = -10283  # was -11492

$ws.Range("H129").Value = 2675.5  # was 2858.3076
$ws.Range("I129").Value = 756  # was 813.125
$ws.Range("K129").Value = 2268  # was 2439.375
$ws.Range("M129").Value = 2732  # was 2560.625

$ws.Range("H134").Value = 4464.3667  # was 4736.857
$ws.Range("I134").Value = 1262.0869  # was 1320.4286
$ws.Range("K134").Value = 3786.2607  # was 3961.2858
$ws.Range("M134").Value = 1283.7393  # was 1108.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 26221.111  # was 26156.37
$ws.Range("J57").Value = 26221.111  # was 26156.37
$ws.Range("L57").Value = 26221.111  # was 26156.37
$ws.Range("N57").Value = -27861.111  # was -27796.37

$ws.Range("H102").Value = 1256.68  # was 1294.75
$ws.Range("I102").Value = 1314.4546  # was 1360.7142
$ws.Range("K102").Value = 1314.4546  # was 1360.7142
$ws.Range("M102").Value = 307.5454  # was 261.2858000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 23393.066  # was 23414.4
$ws.Range("I7").Value = 32975  # was 33007
$ws.Range("K7").Value = 32975  # was 33007
$ws.Range("M7").Value = -32863  # was -32895

$ws.Range("H93").Value = 19812.945  # was 18980.684
$ws.Range("I93").Value = 1300  # was 1458.8235
$ws.Range("K93").Value = 1300  # was 1458.8235
$ws.Range("M93").Value = -52  # was -210.8235

$ws.Range("H100").Value = 27059.8  # was 49092.25
$ws.Range("I100").Value = 4751.25  # was 10801.667
$ws.Range("J100").Value = 52555.285  # was 72066.60000000001
$ws.Range("K100").Value = 4751.25  # was 10801.667
$ws.Range("L100").Value = 52555.285  # was 72066.60000000001
$ws.Range("M100").Value = -4210.25  # was -10260.667
$ws.Range("N100").Value = -53637.285  # was -73148.60000000001

$ws.Range("H122").Value = 16537.125  # was 14119.1
$ws.Range("J122").Value = 4899.5  # was 4786.375
$ws.Range("L122").Value = 14698.5  # was 14359.125
$ws.Range("N122").Value = -19598.5  # was -19259.125

$ws.Range("H126").Value = 23393.066  # was 23414.4
$ws.Range("I126").Value = 32975  # was 33007
$ws.Range("K126").Value = 98925  # was 99021
$ws.Range("M126").Value = -96455  # was -96551

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 29495  # was 28029
$ws.Range("J37").Value = 29495  # was 28029
$ws.Range("L37").Value = 29495  # was 28029
$ws.Range("N37").Value = -29901  # was -28435

$ws.Range("H38").Value = 5000  # was 0
$ws.Range("I38").Value = 5000  # was 0
$ws.Range("K38").Value = 5000  # was 0
$ws.Range("M38").Value = -4527  # newly populated cell

$ws.Range("H39").Value = 17747.5  # was 0
$ws.Range("I39").Value = 5000  # was 0
$ws.Range("J39").Value = 30495  # was 0
$ws.Range("K39").Value = 5000  # was 0
$ws.Range("L39").Value = 30495  # was 0
$ws.Range("M39").Value = -4587  # newly populated cell
$ws.Range("N39").Value = -31321  # newly populated cell

$ws.Range("H54").Value = 46165.332  # was 49999
$ws.Range("J54").Value = 46165.332  # was 49999
$ws.Range("L54").Value = 46165.332  # was 49999
$ws.Range("N54").Value = -47205.332  # was -51039

